# Apply the "previous gates demo" edit to the visit.xlsx workbook.
#
# Summary of the real content changes (see commit message / diff):
#  - survey!F2:   "plot_name"      -> "plot_id"
#  - queries!E2:  "plot_id >= ?"   -> "_id >= ?"
#  - settings:    add a new row (row 5) with setting_name="table_id", value="visit"
#  - selection / active-sheet bookkeeping is updated to match the final
#    state captured in the diff (queries no longer the active tab,
#    settings becomes the active tab, and the selections on survey /
#    queries / settings move to F8 / E3 / B6 respectively).

$wb = $excel.ActiveWorkbook

# --- survey sheet -----------------------------------------------------
$survey = $wb.Worksheets.Item("survey")
$survey.Range("F2").Value = "plot_id"
[void]$survey.Range("F8").Select()

# --- queries sheet ------------------------------------------------------
$queries = $wb.Worksheets.Item("queries")
$queries.Range("E2").Value = "_id >= ?"
[void]$queries.Range("E3").Select()

# --- settings sheet -----------------------------------------------------
$settings = $wb.Worksheets.Item("settings")
$settings.Range("A5").Value = "table_id"
$settings.Range("B5").Value = "visit"

# Activate settings last so it becomes the active/selected tab, and pick
# the final selected cell there, matching the target workbook state.
[void]$settings.Activate()
[void]$settings.Range("B6").Select()
